# Zeitaufzeichnung.xlsx update
# - Update "Projektmanagement" actual hours (H2) 19 -> 15
# - Update "Projektcontroling" actual hours (H5) 12 -> 8
# - Update "Lastenheft" planned time (F10) 1.5 -> 2
# - Move the active selection to L11 (last worked-on cell before saving)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 15
$ws.Range("H5").Value = 8
$ws.Range("F10").Value = 2

$ws.Range("L11").Select()
